$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (the single data record) ---
$ws.Range("B2").Value = "SEI-260007/005119/2024"
$ws.Range("C2").Value = "OPMES"
$ws.Range("F2").Value = "'18"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "Robson"

# "Criação Proc./Pesquisa" (H2) is cleared for this record (value + format).
$ws.Range("H2").Clear()

# "Nº Pesquisa" (I2) / "Finalização da Pesquisa" (J2) get new values.
$ws.Range("I2").Value = "8046/2024"
$ws.Range("J2").Value = 45565

$ws.Range("O2").Value = "Janeiro"
$ws.Range("P2").Value = "'8"
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Value = 0

# Remove the "Observação" column (R) entirely - header and data.
$ws.Range("R1:R2").EntireColumn.Delete()
